$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Price (column D) values for the rows that changed.
$priceUpdates = @{
    2  = "267.72"
    3  = "22.91"
    4  = "6.313"
    5  = "0.06179"
    6  = "3.591"
    7  = "6.698"
    8  = "1.382"
    9  = "0.8335"
    10 = "0.01364"
    12 = "0.08216"
    13 = "0.03415"
    14 = "0.03254"
    15 = "0.09283"
    17 = "0.001720"
    18 = "0.04862"
    19 = "0.006300"
    20 = "0.005367"
    21 = "0.001091"
    22 = "0.0001501"
    23 = "3.761"
    24 = "2.350"
    25 = "0.3348"
    27 = "0.0002684"
    41 = "0.006976"
    42 = "0.1152"
    43 = "0.003330"
    44 = "0.01236"
    45 = "0.00006233"
    47 = "0.7005"
    48 = "0.1793"
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = "'" + $priceUpdates[$row]
}

# Hora (column G) goes from 22 to 23 for every data row (2-51).
for ($row = 2; $row -le 51; $row++) {
    $ws.Range("G$row").NumberFormat = "@"
    $ws.Range("G$row").Value = "'23"
}
